$d = $word.ActiveDocument

$d.Content.Find.Execute("ROLL NO:737819ECR079", $true, $false, $false, $false, $false, $true, 1, $false, "ROLL NO:737819ECR104", 2)
$d.Content.Find.Execute("void setup()", $true, $false, $false, $false, $false, $true, 1, $false, "void setup()", 2)
$d.Content.Find.Execute("pinMode(13,OUTPUT);", $true, $false, $false, $false, $false, $true, 1, $false, "pinMode(13,OUTPUT);", 2)
$d.Content.Find.Execute("pinMode(12,INPUT);", $true, $false, $false, $false, $false, $true, 1, $false, "pinMode(12,INPUT);", 2)
$d.Content.Find.Execute("void loop()", $true, $false, $false, $false, $false, $true, 1, $false, "void loop()", 2)
$d.Content.Find.Execute("pirsensor = digitalRead(2);", $true, $false, $false, $false, $false, $true, 1, $false, "pirsensor = digitalRead(2);", 2)
$d.Content.Find.Execute("digitalWrite(13,HIGH);", $true, $false, $false, $false, $false, $true, 1, $false, "digitalWrite(13,HIGH);", 2)
$d.Content.Find.Execute("delay(10000);", $true, $false, $false, $false, $false, $true, 1, $false, "delay(10000);", 2)
$d.Content.Find.Execute("digitalWrite(13,LOW);", $true, $false, $false, $false, $false, $true, 1, $false, "digitalWrite(13,LOW);", 2)
$d.Content.Find.Execute("delay(120);", $true, $false, $false, $false, $false, $true, 1, $false, "delay(120);", 2)
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute(" SMART ", $true, $false, $false, $false, $false, $true, 1, $false, " SMART ", 2)
$hdr.Range.Find.Execute(" HOME AUTOMATION ASSIGNMENT 1", $true, $false, $false, $false, $false, $true, 1, $false, " HOME AUTOMATION ASSIGNMENT 1", 2)
